$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new entry row for "photbiology"
$ws.Range("A58").Value = "photbiology"
$ws.Range("B58").Value = "https://www.r4photobiology.info/"
$ws.Range("C58").Value = "https://bitbucket.org/aphalo/photobiology/src/master/"
$ws.Range("D58").Value = "https://bulletin.uv4plants.org/index.php/uv4pbulletin/article/view/16"
$ws.Range("E58").Value = "R"
$ws.Range("F58").Value = "Photobiology workflow tools"
$ws.Range("G58").Value = "UV-Vis"

# Fix capitalization of "Mid-Ir" -> "Mid-IR" for the geoSpectral row (row 57)
$ws.Range("G57").Value = "Mid-IR"

# Update the view: scroll/selection state seen in the saved workbook
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Range("G58").Select()
